$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report date
$ws.Range("C2").Value = "July 25, 2024"

# Center-align the header/value cells (A1:B2 stay regular font, C1:D2 becomes bold)
$ws.Range("A1:D2").HorizontalAlignment = -4108

# Make the Date column (C:D, merged) bold
$ws.Range("C1:D2").Font.Bold = $true

# Merge the Date header/value cells across C:D
$ws.Range("C1:D1").Merge()
$ws.Range("C2:D2").Merge()

# Widen column A to fit the longer title text
$ws.Columns("A").ColumnWidth = 19.86

# Match the author's final selection
$ws.Range("D4").Select() | Out-Null
